$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-12 (fiscal years 2017-2027) with refreshed
# Controller (B) and Five Year Plan (C) figures for the FY24-FY28 plan.
$ws.Range("B2").Value = 31934098.00767143
$ws.Range("C2").Value = 31934098.00767143

$ws.Range("B3").Value = 78054304.06804928
$ws.Range("C3").Value = 78054304.06804928

$ws.Range("B4").Value = 77406811.51657739
$ws.Range("C4").Value = 77406811.51657739

$ws.Range("B5").Value = 70659178.24173306
$ws.Range("C5").Value = 70659178.24173306

$ws.Range("B6").Value = 69581516.35580111
$ws.Range("C6").Value = 69581516.35580111

$ws.Range("B7").Value = 75427039.33037321
$ws.Range("C7").Value = 75427039.33037321

$ws.Range("B8").Value = 68087000
$ws.Range("C8").Value = 68087000

$ws.Range("B9").Value = 74167236.31064345
$ws.Range("C9").Value = 68938000

$ws.Range("B10").Value = 73325023.28403327
$ws.Range("C10").Value = 69751000

$ws.Range("B11").Value = 72677585.91664648
$ws.Range("C11").Value = 69053000

$ws.Range("B12").Value = 72029638.8409044
$ws.Range("C12").Value = 68680000

# Add the new FY2028 row, matching the formatting already used for the
# "fiscal_year" column (copy format from A2 so it picks up the same
# bold / centered / bordered style instead of minting a new one).
$ws.Range("A13").Value = 2028
$ws.Range("A2").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B13").Value = 71539973.72882061
$ws.Range("C13").Value = 68385000
